$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6").Value = 30
$ws.Range("B6").Value = "first 30 in data folder"
$ws.Range("C6").Value = 95.7

$ws.Range("A7").Value = 44
$ws.Range("B7").Value = "all"
